# Generate Report for Handoff
#
# The localization-status report was regenerated; the only functional
# change is the "Latest Handoff Datetime" value recorded for the
# c960c4b1-dc6e-4e55-8edd-cff68a2aeb04 file's zh-cn handoff xliff, which
# moved from 2016-08-19 14:47:52 to 2016-08-19 14:48:17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

# Row 5 of the zh-cn table is the c960c4b1-dc6e-4e55-8edd-cff68a2aeb04 entry;
# column H is "Latest Handoff Datetime" (stored as plain text, matching the
# existing yyyy-mm-dd HH:mm:ss-formatted text cells in this column).
$ws.Range("H5").Value = "2016-08-19 14:48:17"
